$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 243.33333
$ws.Range("I8").Value = 243.33333
$ws.Range("K8").Value = 729.99999
$ws.Range("M8").Value = -590.99999
$ws.Range("H43").Value = 1100.1538
$ws.Range("I43").Value = 800
$ws.Range("J43").Value = 1233.5555
$ws.Range("K43").Value = 800
$ws.Range("L43").Value = 1233.5555
$ws.Range("M43").Value = -731
$ws.Range("N43").Value = -1371.5555
$ws.Range("H116").Value = 7399.4736
$ws.Range("I116").Value = 10808.182
$ws.Range("K116").Value = 10808.182
$ws.Range("M116").Value = -7366.182000000001
$ws.Range("H132").Value = 1089.742
$ws.Range("I132").Value = 888.8276
$ws.Range("K132").Value = 2666.4828
$ws.Range("M132").Value = -136.4827999999998

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5012.4243
$ws.Range("I32").Value = 4005.2932
$ws.Range("K32").Value = 4005.2932
$ws.Range("M32").Value = -3718.2932
$ws.Range("H74").Value = 2282
$ws.Range("I74").Value = 2150.2666
$ws.Range("J74").Value = 2501.5557
$ws.Range("K74").Value = 2150.2666
$ws.Range("L74").Value = 2501.5557
$ws.Range("M74").Value = -1276.2666
$ws.Range("N74").Value = -4249.5557
$ws.Range("H77").Value = 2282
$ws.Range("I77").Value = 2150.2666
$ws.Range("J77").Value = 2501.5557
$ws.Range("K77").Value = 10751.333
$ws.Range("L77").Value = 12507.7785
$ws.Range("M77").Value = -6383.332999999999
$ws.Range("N77").Value = -21243.7785
$ws.Range("H88").Value = 2432.2222
$ws.Range("I88").Value = 2269.4285
$ws.Range("K88").Value = 2269.4285
$ws.Range("M88").Value = -1863.4285
$ws.Range("H91").Value = 2432.2222
$ws.Range("I91").Value = 2269.4285
$ws.Range("K91").Value = 2269.4285
$ws.Range("M91").Value = -865.4285

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 207.90475
$ws.Range("I80").Value = 68.85714
$ws.Range("J80").Value = 277.42856
$ws.Range("K80").Value = 68.85714
$ws.Range("L80").Value = 277.42856
$ws.Range("M80").Value = 929.14286
$ws.Range("N80").Value = -2273.42856
$ws.Range("H83").Value = 207.90475
$ws.Range("I83").Value = 68.85714
$ws.Range("J83").Value = 277.42856
$ws.Range("K83").Value = 344.2857
$ws.Range("L83").Value = 1387.1428
$ws.Range("M83").Value = 4647.7143
$ws.Range("N83").Value = -11371.1428
$ws.Range("H99").Value = 111112424
$ws.Range("I99").Value = 125001100
$ws.Range("K99").Value = 125001100
$ws.Range("M99").Value = -124999602

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1543.2222
$ws.Range("J16").Value = 1481.6666
$ws.Range("L16").Value = 1481.6666
$ws.Range("N16").Value = -2055.6666
$ws.Range("H99").Value = 3459.5557
$ws.Range("I99").Value = 1962
$ws.Range("J99").Value = 5812.857
$ws.Range("K99").Value = 1962
$ws.Range("L99").Value = 5812.857
$ws.Range("M99").Value = -464
$ws.Range("N99").Value = -8808.857
$ws.Range("H113").Value = 1543.2222
$ws.Range("J113").Value = 1481.6666
$ws.Range("L113").Value = 1481.6666
$ws.Range("N113").Value = -5821.6666
$ws.Range("H122").Value = 1489.9286
$ws.Range("I122").Value = 1329.1111
$ws.Range("J122").Value = 1779.4
$ws.Range("K122").Value = 3987.3333
$ws.Range("L122").Value = 5338.200000000001
$ws.Range("M122").Value = -1537.3333
$ws.Range("N122").Value = -10238.2
$ws.Range("H126").Value = 3459.5557
$ws.Range("I126").Value = 1962
$ws.Range("J126").Value = 5812.857
$ws.Range("K126").Value = 5886
$ws.Range("L126").Value = 17438.571
$ws.Range("M126").Value = -3416
$ws.Range("N126").Value = -22378.571

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 23.666666
$ws.Range("J2").Value = 99
$ws.Range("L2").Value = 594
$ws.Range("N2").Value = -820
$ws.Range("H5").Value = 300524.9
$ws.Range("I5").Value = 766.3333
$ws.Range("K5").Value = 2298.9999
$ws.Range("M5").Value = -2186.9999
$ws.Range("H68").Value = 491
$ws.Range("I68").Value = 428.14285
$ws.Range("J68").Value = 711
$ws.Range("K68").Value = 1284.42855
$ws.Range("L68").Value = 2133
$ws.Range("M68").Value = -473.4285500000001
$ws.Range("N68").Value = -3755
$ws.Range("H71").Value = 491
$ws.Range("I71").Value = 428.14285
$ws.Range("J71").Value = 711
$ws.Range("K71").Value = 3853.28565
$ws.Range("L71").Value = 6399
$ws.Range("M71").Value = 202.7143499999997
$ws.Range("N71").Value = -14511
$ws.Range("H122").Value = 780.8
$ws.Range("J122").Value = 977.5
$ws.Range("L122").Value = 8797.5
$ws.Range("N122").Value = -13697.5
$ws.Range("H131").Value = 1786688.5
$ws.Range("J131").Value = 1055.7142
$ws.Range("L131").Value = 3167.1426
$ws.Range("N131").Value = -13247.1426
$ws.Range("H132").Value = 10104292
$ws.Range("I132").Value = 833.3333
$ws.Range("J132").Value = 13893090
$ws.Range("K132").Value = 7499.9997
$ws.Range("L132").Value = 125037810
$ws.Range("M132").Value = -4969.9997
$ws.Range("N132").Value = -125042870
$ws.Range("H135").Value = 300524.9
$ws.Range("I135").Value = 766.3333
$ws.Range("K135").Value = 6896.9997
$ws.Range("M135").Value = -4361.9997

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 47004.5
$ws.Range("J6").Value = 47004.5
$ws.Range("L6").Value = 47004.5
$ws.Range("N6").Value = -47230.5
$ws.Range("H16").Value = 47004.5
$ws.Range("J16").Value = 47004.5
$ws.Range("L16").Value = 47004.5
$ws.Range("N16").Value = -47504.5
$ws.Range("H113").Value = 71429810
$ws.Range("I113").Value = 125000780
$ws.Range("J113").Value = 1850
$ws.Range("K113").Value = 125000780
$ws.Range("L113").Value = 1850
$ws.Range("M113").Value = -124998610
$ws.Range("N113").Value = -6190
$ws.Range("H132").Value = 5953.7334
$ws.Range("J132").Value = 3062
$ws.Range("L132").Value = 9186
$ws.Range("N132").Value = -14246

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 34357.902
$ws.Range("I7").Value = 47276.816
$ws.Range("J7").Value = 2778.3333
$ws.Range("K7").Value = 47276.816
$ws.Range("L7").Value = 2778.3333
$ws.Range("M7").Value = -47164.816
$ws.Range("N7").Value = -3002.3333
$ws.Range("H126").Value = 34357.902
$ws.Range("I126").Value = 47276.816
$ws.Range("J126").Value = 2778.3333
$ws.Range("K126").Value = 141830.448
$ws.Range("L126").Value = 8334.999899999999
$ws.Range("M126").Value = -139360.448
$ws.Range("N126").Value = -13274.9999
$ws.Range("H136").Value = 5762.1177
$ws.Range("I136").Value = 6740.28
$ws.Range("K136").Value = 20220.84
$ws.Range("M136").Value = -17670.84

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1715.075
$ws.Range("I132").Value = 1317.9546
$ws.Range("J132").Value = 2200.4443
$ws.Range("K132").Value = 3953.8638
$ws.Range("L132").Value = 6601.3329
$ws.Range("M132").Value = -1423.8638
$ws.Range("N132").Value = -11661.3329
